# Replace the Python-dict-style question dump with pretty-printed JSON,
# move it from A2 up to A1 (taking on the default/unstyled look),
# and drop the old bold/bordered A1 placeholder cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$questionsText = @'
questions = [
    {
        "title": "You are a Firebase developer working on an e-commerce store with a Firebase backend. You want to secure the storage of user passwords on your site.What should you do?",
        "ques_type": 2,
        "options": [
            "Use a salted hash algorithm to hash passwords securely.",
            "Use a symmetric encryption algorithm to store password digests.",
            "Use a database with restricted access permissions to store passwords.",
            "Use a secure token-based authentication mechanism for password storage."
        ],
        "score": "Use a salted hash algorithm to hash passwords securely."
    },
    {
        "title": "You are a Firebase developer working on a mobile app, and you are using Firestore for your backend database. You have implemented security rules to restrict access to sensitive data. However, you notice that a user with unauthorized access is able to read and modify that sensitive data.What should you do to address this vulnerability?",
        "ques_type": 2,
        "options": [
            "Review and update Firestore security rules for proper access restrictions on sensitive data.",
            "Implement Firebase Authentication to enforce access controls and authenticate users.",
            "Utilize Firebase Cloud Functions for server-side validation and authorization checks on sensitive data.",
            "Implement additional security measures to strengthen the protection of sensitive data."
        ],
        "score": "Review and update Firestore security rules for proper access restrictions on sensitive data."
    },
    {
        "title": "You are a developer working with Firebase Cloud Functions. You encounter an unexpected error during a critical database operation. The function crashes, but the error message and details are not logged. This lack of information poses a challenge in identifying the error's cause and troubleshooting. What should you do?",
        "ques_type": 2,
        "options": [
            "Enable Firebase Crashlytics to capture and report errors within the Cloud Function.",
            "Utilize Firebase Performance Monitoring to track the Cloud Function's performance.",
            "Implement Firebase Remote Config for real-time adjustment of the Cloud Function's behavior.",
            "Integrate Firebase Analytics to gain insights into user interactions within the Cloud Function."
        ],
        "score": "Enable Firebase Crashlytics to capture and report errors within the Cloud Function."
    },
    {
        "title": "You are a Firebase developer configuring Message Topics for a messaging app. You find that, despite correct topic subscriptions, some users experience a five-minute delay in receiving notifications.What should you do to troubleshoot this delayed notification issue?",
        "ques_type": 2,
        "options": [
            "Verify the accuracy of user topic subscriptions and device tokens.",
            "Increase the notification priority level for all users to expedite delivery.",
            "Upgrade the server hardware to handle higher notification load.",
            "Modify the authentication flow to ensure notifications are sent only to verified users."
        ],
        "score": "Verify the accuracy of user topic subscriptions and device tokens."
    }
]
'@

# Remove row 1 (the bold/bordered "0" placeholder cell). This shifts
# row 2 (the shared-string question dump, A2) up to row 1 with no style,
# matching the diff: A1/s=1 cell is gone, A2 content now lives at A1.
$ws.Rows(1).Delete()

# Update the text to the reformatted (pretty JSON) version.
$ws.Range("A1").Value = $questionsText

